$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Find last used row based on column A (Beteckning)
$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row  # xlUp = -4162

# Column C holds "Förändrad" (changed) date serials; bump every value by 1 day
for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 3)
    $current = $cell.Value2
    if ($current -ne $null) {
        $cell.Value2 = $current + 1
    }
}
